# Updated OVW demo sheet for Architechture pages
#
# Every locale worksheet has a small A1:D2/A1:D3 table whose column C holds
# "collaboration-architecture-var1/3/4" labels (one per row). This drops the
# stale "collaboration-" prefix so the labels read "architecture-var1/3/4",
# across every worksheet in the workbook.

$wb = $excel.ActiveWorkbook

$oldPrefix = "collaboration-architecture-var"
$newPrefix = "architecture-var"

foreach ($ws in $wb.Worksheets) {
    for ($r = 1; $r -le 3; $r++) {
        $cell = $ws.Cells.Item($r, 3)
        $val = $cell.Value2
        if ($val -ne $null) {
            $text = $val.ToString()
            if ($text.StartsWith($oldPrefix)) {
                $suffix = $text.Substring($oldPrefix.Length)
                $cell.Value = $newPrefix + $suffix
            }
        }
    }
}

# The last locale tab (uk_ua) ends up the active sheet/tab in the saved file.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$lastSheet.Activate()
